# Add the new "Prueba Bautista" registration row beneath the existing data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Prueba"
$ws.Range("B3").Value = "Bautista"
$ws.Range("C3").Value = "jbautmqpea@gmail.com"
$ws.Range("D3").Value = "1730224973880-ATS CV harvard (1).pdf"
